$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 with new TPM-derived values
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.02640533333333334
$ws.Range("H2").Value = 0.07921600000000001
$ws.Range("M2").Value = 29.75868033333333
$ws.Range("N2").Value = 89.27604099999999
$ws.Range("O2").Value = 0.4948552779010537
$ws.Range("P2").Value = 0.4948552779010535
$ws.Range("Q2").Value = 0.7857878737617778
$ws.Range("R2").Value = 7.072090863856
$ws.Range("S2").Value = 0.4948552779010537
$ws.Range("T2").Value = 0.4948552779010535

$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.02640533333333334
$ws.Range("H3").Value = 0.07921600000000001
$ws.Range("M3").Value = 17.55525033333333
$ws.Range("N3").Value = 52.665751
$ws.Range("O3").Value = 0.2919251856942525
$ws.Range("P3").Value = 0.2919251856942524
$ws.Range("Q3").Value = 0.4635522368017778
$ws.Range("R3").Value = 4.171970131216001
$ws.Range("S3").Value = 0.2919251856942525
$ws.Range("T3").Value = 0.2919251856942524

$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.02640533333333334
$ws.Range("H4").Value = 0.07921600000000001
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1421396666666667
$ws.Range("N4").Value = 0.426419
$ws.Range("O4").Value = 0.002363631836533717
$ws.Range("P4").Value = 0.002363631836533717
$ws.Range("Q4").Value = 0.003753245278222223
$ws.Range("R4").Value = 0.033779207504
$ws.Range("S4").Value = 0.002363631836533717
$ws.Range("T4").Value = 0.002363631836533717

# Add new rows 5 and 6 for additional target clusters
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ndp"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.02640533333333334
$ws.Range("H5").Value = 0.07921600000000001
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.42872866666667
$ws.Range("N5").Value = 37.286186
$ws.Range("O5").Value = 0.2066765699758167
$ws.Range("P5").Value = 0.2066765699758166
$ws.Range("Q5").Value = 0.3281847233528889
$ws.Range("R5").Value = 2.953662510176001
$ws.Range("S5").Value = 0.2066765699758167
$ws.Range("T5").Value = 0.2066765699758166

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ndp"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.02640533333333334
$ws.Range("H6").Value = 0.07921600000000001
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.251329
$ws.Range("N6").Value = 0.753987
$ws.Range("O6").Value = 0.004179334592343558
$ws.Range("P6").Value = 0.004179334592343557
$ws.Range("Q6").Value = 0.006636426021333333
$ws.Range("R6").Value = 0.05972783419200001
$ws.Range("S6").Value = 0.004179334592343558
$ws.Range("T6").Value = 0.004179334592343557

